$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Description" note element previously had a type="local" attribute;
# it has been removed. Update cell V2 accordingly, keeping its style.
$ws.Range("V2").Value = '<mods:note displayLabel="Description">'

# Reflect the active selection that results from editing this cell.
$ws.Range("V2").Select()
